# Insert a new weekly data row for "Poroto granado" (Terminal La Palmera de
# La Serena) just below the header/first-data block, at row 35, pushing the
# existing rows 35-75 down by one (to 36-76).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 35:75 down to 36:76 to make room for the new entry at row 35.
$ws.Rows("35:35").Insert()

# Populate the newly inserted row 35 with this week's price observation.
$ws.Range("A35").Value = 8
$ws.Range("B35").Value = "Terminal La Palmera de La Serena"
$ws.Range("C35").Value = "Coquimbo"
$ws.Range("D35").Value = (Get-Date -Year 2022 -Month 3 -Day 9 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("D35").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E35").Value = 4
$ws.Range("F35").Value = 100112030
$ws.Range("G35").Value = "Poroto granado"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 400
$ws.Range("K35").Value = 29000
$ws.Range("L35").Value = 30000
$ws.Range("M35").Value = 29500
$ws.Range("N35").Value = '$/malla 25 kilos'
$ws.Range("O35").Value = "Provincia del Elquí"
$ws.Range("P35").Value = 1180
$ws.Range("Q35").Value = 25
$ws.Range("R35").Value = "Hortaliza"
